$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed TPM-derived values (ligand/receptor expression,
# specificity scores, and edge weights) per commit "update scripts wuth new tpm".

# Row 2
$ws.Range("G2").Value = 2.269101333333333
$ws.Range("H2").Value = 6.807304
$ws.Range("I2").Value = 0.02891211995713196
$ws.Range("J2").Value = 0.02891211995713196
$ws.Range("M2").Value = 14.861848
$ws.Range("N2").Value = 44.585544
$ws.Range("O2").Value = 0.09055189482833943
$ws.Range("P2").Value = 0.09055189482833945
$ws.Range("Q2").Value = 33.72303911259733
$ws.Range("R2").Value = 303.507352013376
$ws.Range("S2").Value = 0.002618047245622547
$ws.Range("T2").Value = 0.002618047245622547

# Row 3
$ws.Range("G3").Value = 2.269101333333333
$ws.Range("H3").Value = 6.807304
$ws.Range("I3").Value = 0.02891211995713196
$ws.Range("J3").Value = 0.02891211995713196
$ws.Range("O3").Value = 0.1893562842131466
$ws.Range("P3").Value = 0.1893562842131466
$ws.Range("Q3").Value = 70.51944512968444
$ws.Range("R3").Value = 634.6750061671601
$ws.Range("S3").Value = 0.005474691603807266
$ws.Range("T3").Value = 0.005474691603807267

# Row 4
$ws.Range("G4").Value = 2.269101333333333
$ws.Range("H4").Value = 6.807304
$ws.Range("I4").Value = 0.02891211995713196
$ws.Range("J4").Value = 0.02891211995713196
$ws.Range("M4").Value = 18.10188466666667
$ws.Range("N4").Value = 54.305654
$ws.Range("O4").Value = 0.1102931450066459
$ws.Range("P4").Value = 0.1102931450066459
$ws.Range("Q4").Value = 41.07501063297956
$ws.Range("R4").Value = 369.675095696816
$ws.Range("S4").Value = 0.003188808638881495
$ws.Range("T4").Value = 0.003188808638881496

# Row 5
$ws.Range("G5").Value = 2.269101333333333
$ws.Range("H5").Value = 6.807304
$ws.Range("I5").Value = 0.02891211995713196
$ws.Range("J5").Value = 0.02891211995713196
$ws.Range("M5").Value = 100.0833306666667
$ws.Range("N5").Value = 300.249992
$ws.Range("O5").Value = 0.609798675951868
$ws.Range("P5").Value = 0.6097986759518681
$ws.Range("Q5").Value = 227.0992190601742
$ws.Range("R5").Value = 2043.892971541568
$ws.Range("S5").Value = 0.01763057246882065
$ws.Range("T5").Value = 0.01763057246882065

# Row 6
$ws.Range("I6").Value = 0.7238963226334669
$ws.Range("J6").Value = 0.7238963226334669
$ws.Range("M6").Value = 14.861848
$ws.Range("N6").Value = 44.585544
$ws.Range("O6").Value = 0.09055189482833943
$ws.Range("P6").Value = 0.09055189482833945
$ws.Range("Q6").Value = 844.3512284062693
$ws.Range("R6").Value = 7599.161055656424
$ws.Range("S6").Value = 0.06555018367372736
$ws.Range("T6").Value = 0.06555018367372738

# Row 7
$ws.Range("I7").Value = 0.7238963226334669
$ws.Range("J7").Value = 0.7238963226334669
$ws.Range("O7").Value = 0.1893562842131466
$ws.Range("P7").Value = 0.1893562842131466
$ws.Range("S7").Value = 0.1370743178094344
$ws.Range("T7").Value = 0.1370743178094344

# Row 8
$ws.Range("I8").Value = 0.7238963226334669
$ws.Range("J8").Value = 0.7238963226334669
$ws.Range("M8").Value = 18.10188466666667
$ws.Range("N8").Value = 54.305654
$ws.Range("O8").Value = 0.1102931450066459
$ws.Range("P8").Value = 0.1102931450066459
$ws.Range("Q8").Value = 1028.428534242082
$ws.Range("R8").Value = 9255.856808178734
$ws.Range("S8").Value = 0.07984080208199068
$ws.Range("T8").Value = 0.07984080208199071

# Row 9
$ws.Range("I9").Value = 0.7238963226334669
$ws.Range("J9").Value = 0.7238963226334669
$ws.Range("M9").Value = 100.0833306666667
$ws.Range("N9").Value = 300.249992
$ws.Range("O9").Value = 0.609798675951868
$ws.Range("P9").Value = 0.6097986759518681
$ws.Range("Q9").Value = 5686.06832685887
$ws.Range("R9").Value = 51174.61494172983
$ws.Range("S9").Value = 0.4414310190683144
$ws.Range("T9").Value = 0.4414310190683144

# Row 10
$ws.Range("G10").Value = 18.57257166666666
$ws.Range("H10").Value = 55.717715
$ws.Range("I10").Value = 0.2366454120188096
$ws.Range("J10").Value = 0.2366454120188096
$ws.Range("M10").Value = 14.861848
$ws.Range("N10").Value = 44.585544
$ws.Range("O10").Value = 0.09055189482833943
$ws.Range("P10").Value = 0.09055189482833945
$ws.Range("Q10").Value = 276.0227370791067
$ws.Range("R10").Value = 2484.20463371196
$ws.Range("S10").Value = 0.0214286904607363
$ws.Range("T10").Value = 0.0214286904607363

# Row 11
$ws.Range("G11").Value = 18.57257166666666
$ws.Range("H11").Value = 55.717715
$ws.Range("I11").Value = 0.2366454120188096
$ws.Range("J11").Value = 0.2366454120188096
$ws.Range("O11").Value = 0.1893562842131466
$ws.Range("P11").Value = 0.1893562842131466
$ws.Range("Q11").Value = 577.2009514624139
$ws.Range("R11").Value = 5194.808563161725
$ws.Range("S11").Value = 0.04481029589597087
$ws.Range("T11").Value = 0.04481029589597089

# Row 12
$ws.Range("G12").Value = 18.57257166666666
$ws.Range("H12").Value = 55.717715
$ws.Range("I12").Value = 0.2366454120188096
$ws.Range("J12").Value = 0.2366454120188096
$ws.Range("M12").Value = 18.10188466666667
$ws.Range("N12").Value = 54.305654
$ws.Range("O12").Value = 0.1102931450066459
$ws.Range("P12").Value = 0.1102931450066459
$ws.Range("Q12").Value = 336.1985502734011
$ws.Range("R12").Value = 3025.78695246061
$ws.Range("S12").Value = 0.02610036674294803
$ws.Range("T12").Value = 0.02610036674294803

# Row 13
$ws.Range("G13").Value = 18.57257166666666
$ws.Range("H13").Value = 55.717715
$ws.Range("I13").Value = 0.2366454120188096
$ws.Range("J13").Value = 0.2366454120188096
$ws.Range("M13").Value = 100.0833306666667
$ws.Range("N13").Value = 300.249992
$ws.Range("O13").Value = 0.609798675951868
$ws.Range("P13").Value = 0.6097986759518681
$ws.Range("Q13").Value = 1858.804831445364
$ws.Range("R13").Value = 16729.24348300828
$ws.Range("S13").Value = 0.1443060589191544
$ws.Range("T13").Value = 0.1443060589191544

# Row 14
$ws.Range("G14").Value = 0.8276899999999999
$ws.Range("H14").Value = 2.48307
$ws.Range("I14").Value = 0.01054614539059158
$ws.Range("J14").Value = 0.01054614539059158
$ws.Range("M14").Value = 14.861848
$ws.Range("N14").Value = 44.585544
$ws.Range("O14").Value = 0.09055189482833943
$ws.Range("P14").Value = 0.09055189482833945
$ws.Range("Q14").Value = 12.30100297112
$ws.Range("R14").Value = 110.70902674008
$ws.Range("S14").Value = 0.0009549734482532256
$ws.Range("T14").Value = 0.0009549734482532257

# Row 15
$ws.Range("G15").Value = 0.8276899999999999
$ws.Range("H15").Value = 2.48307
$ws.Range("I15").Value = 0.01054614539059158
$ws.Range("J15").Value = 0.01054614539059158
$ws.Range("O15").Value = 0.1893562842131466
$ws.Range("P15").Value = 0.1893562842131466
$ws.Range("Q15").Value = 25.72306431711666
$ws.Range("R15").Value = 231.50757885405
$ws.Range("S15").Value = 0.001996978903934025
$ws.Range("T15").Value = 0.001996978903934025

# Row 16
$ws.Range("G16").Value = 0.8276899999999999
$ws.Range("H16").Value = 2.48307
$ws.Range("I16").Value = 0.01054614539059158
$ws.Range("J16").Value = 0.01054614539059158
$ws.Range("M16").Value = 18.10188466666667
$ws.Range("N16").Value = 54.305654
$ws.Range("O16").Value = 0.1102931450066459
$ws.Range("P16").Value = 0.1102931450066459
$ws.Range("Q16").Value = 14.98274891975333
$ws.Range("R16").Value = 134.84474027778
$ws.Range("S16").Value = 0.001163167542825688
$ws.Range("T16").Value = 0.001163167542825688

# Row 17
$ws.Range("G17").Value = 0.8276899999999999
$ws.Range("H17").Value = 2.48307
$ws.Range("I17").Value = 0.01054614539059158
$ws.Range("J17").Value = 0.01054614539059158
$ws.Range("M17").Value = 100.0833306666667
$ws.Range("N17").Value = 300.249992
$ws.Range("O17").Value = 0.609798675951868
$ws.Range("P17").Value = 0.6097986759518681
$ws.Range("Q17").Value = 82.83797195949332
$ws.Range("R17").Value = 745.5417476354398
$ws.Range("S17").Value = 0.006431025495578642
$ws.Range("T17").Value = 0.006431025495578643
